$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("modelos por sujeto sano")

$ws1 = $wb.Worksheets.Item("full_signals - without decay")
$ws2 = $wb.Worksheets.Item("full_signals - with decay")

# Row 6 - unet_model_vscdv3.keras
$ws.Cells.Item(6, 2).Value = "1_HEMU"
$ws.Cells.Item(6, 3).Value = "unet_model_vscdv3.keras"
$ws.Cells.Item(6, 4).Value = 50
$ws.Cells.Item(6, 5).Value = "80/20"
$ws.Cells.Item(6, 7).Value = 0.0001
$ws.Cells.Item(6, 8).Value = 200
$ws.Cells.Item(6, 9).Value = 4
$ws.Cells.Item(6, 12).Value = "Adam"
$ws.Cells.Item(6, 14).Value = 0.0017
$ws.Cells.Item(6, 15).Value = 0.0042
$ws.Cells.Item(6, 16).Value = 0.1146
$ws.Cells.Item(6, 17).Value = 0.2972
$ws.Cells.Item(6, 19).Value = 87

# Row 7 - unet_model_vscdv4.keras
$ws.Cells.Item(7, 3).Value = "unet_model_vscdv4.keras"
$ws.Cells.Item(7, 4).Value = 50
$ws.Cells.Item(7, 5).Value = "80/20"
$ws.Cells.Item(7, 6).Value = "z-core"
$ws.Cells.Item(7, 7).Value = 0.0001
$ws.Cells.Item(7, 8).Value = 200
$ws.Cells.Item(7, 9).Value = 4
$ws.Cells.Item(7, 12).Value = "Adam"
$ws.Cells.Item(7, 14).Value = 0.1597
$ws.Cells.Item(7, 15).Value = 0.0097
$ws.Cells.Item(7, 16).Value = 4.6406
$ws.Cells.Item(7, 17).Value = 0.2956
$ws.Cells.Item(7, 19).Value = 91

# Row 8 - unet_model_vscdv5.keras
$ws.Cells.Item(8, 3).Value = "unet_model_vscdv5.keras"
$ws.Cells.Item(8, 4).Value = 50
$ws.Cells.Item(8, 5).Value = "70/30"
$ws.Cells.Item(8, 7).Value = 0.01
$ws.Cells.Item(8, 8).Value = 200
$ws.Cells.Item(8, 9).Value = 16
$ws.Cells.Item(8, 12).Value = "Adam"
$ws.Cells.Item(8, 14).Value = 0.0068
$ws.Cells.Item(8, 15).Value = 0.0172
$ws.Cells.Item(8, 16).Value = 0.1318
$ws.Cells.Item(8, 17).Value = 0.3446
$ws.Cells.Item(8, 19).Value = 82

# Row 9 - unet_model_vscdv6.keras
$ws.Cells.Item(9, 3).Value = "unet_model_vscdv6.keras"
$ws.Cells.Item(9, 4).Value = 50
$ws.Cells.Item(9, 5).Value = "70/30"
$ws.Cells.Item(9, 7).Value = 0.01
$ws.Cells.Item(9, 8).Value = 300
$ws.Cells.Item(9, 9).Value = 32
$ws.Cells.Item(9, 12).Value = "adam"

# Update view/selection state on each sheet to match the session's final state
$ws2.Activate()
$ws2.Range("V20").Select()

$ws1.Activate()

$ws.Activate()
$ws.Range("H9").Select()
